# Council workbook update: record the official winner flag for each race
# and correct a few vote-percentage figures that were re-tabulated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected vote percentages (re-tabulated figures)
$ws.Range("E14").Value = 0.88
$ws.Range("E15").Value = 0.9
$ws.Range("E33").Value = 0.91

# Mark the winning candidate (highest vote-getter) in each race as TRUE
$winnerRows = @(6, 7, 9, 11, 14, 15, 17, 21, 22, 26, 30, 32, 33)
foreach ($r in $winnerRows) {
    $ws.Range("G$r").Value = $true
}

# Highlight the full results table now that winners have been added
$ws.Range("A1:G34").Select() | Out-Null
